# Remove the "vehicle_id" column (column A) from the vehicles sheet.
# This shifts reference/year/engine_capacity/brand/price one column to the
# left (B->A, C->B, D->C, E->D, F->E) and drops the old vehicle_id values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete()

# Mirror Excel's behaviour of leaving the freshly-shifted column selected
# after a "Delete Column" operation.
$ws.Columns.Item(1).Select()
